$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''60.987.47'
$ws.Range("E2").Value = '  +7.65%  '
$ws.Range("D3").Value = '''2.677.21'
$ws.Range("E3").Value = '  +11.18%  '
$ws.Range("D4").Value = '''0.997'
$ws.Range("E4").Value = '  -0.59%  '
$ws.Range("D5").Value = '''513.31'
$ws.Range("E5").Value = '  +5.84%  '
$ws.Range("D6").Value = '''157.70'
$ws.Range("E6").Value = '  +3.22%  '
$ws.Range("E7").Value = '  +1.13%  '
$ws.Range("D8").Value = '''0.996'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("D9").Value = '''2.676.27'
$ws.Range("E9").Value = '  +10.44%  '
$ws.Range("D10").Value = '''6.32'
$ws.Range("E10").Value = '  +12.36%  '
$ws.Range("E11").Value = '  +5.85%  '
$ws.Range("D12").Value = '''0.348'
$ws.Range("E12").Value = '  +4.26%  '
$ws.Range("E13").Value = '  +0.97%  '
$ws.Range("D14").Value = '''3.113.17'
$ws.Range("E14").Value = '  +9.90%  '
$ws.Range("D15").Value = '''61.036.09'
$ws.Range("E15").Value = '  +7.00%  '
$ws.Range("D16").Value = '''21.89'
$ws.Range("E16").Value = '  +5.72%  '
$ws.Range("E17").Value = '  +5.47%  '
$ws.Range("D18").Value = '''2.670.90'
$ws.Range("E18").Value = '  +10.03%  '
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").Value = '''350.61'
$ws.Range("E20").Value = '  +8.35%  '
$ws.Range("E21").Value = '  +5.98%  '
$ws.Range("E22").Value = '  +5.33%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").Value = '''60.21'
$ws.Range("E24").Value = '  +3.60%  '
$ws.Range("E25").Value = '  +3.80%  '
$ws.Range("D26").Value = '''2.766.58'
$ws.Range("E26").Value = '  +9.63%  '
$ws.Range("E27").Value = '  +5.21%  '
$ws.Range("D28").Value = '''0.990'
$ws.Range("D29").Value = '''0.0₃0870'
$ws.Range("E29").Value = '  +11.67%  '
$ws.Range("E30").Value = '  +4.59%  '
$ws.Range("D31").Value = '''0.997'
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("E32").Value = '  +5.73%  '
$ws.Range("E33").Value = '  +6.12%  '
$ws.Range("E34").Value = '  +4.21%  '
$ws.Range("D35").Value = '''5.73'
$ws.Range("E35").Value = '  +7.12%  '
$ws.Range("E36").Value = '  +9.57%  '
$ws.Range("E37").Value = '  +6.54%  '
$ws.Range("E38").Value = '  +11.24%  '
$ws.Range("D39").Value = '''0.869'
$ws.Range("E39").Value = '  +2.92%  '
$ws.Range("D40").Value = '''310.13'
$ws.Range("E40").Value = '  +17.53%  '
$ws.Range("E41").Value = '  +7.79%  '
$ws.Range("D42").Value = '''0.836'
$ws.Range("E42").Value = '  +29.94%  '
$ws.Range("D43").Value = '''35.42'
$ws.Range("E43").Value = '  +3.93%  '
$ws.Range("D44").Value = '''0.648'
$ws.Range("E44").Value = '  +9.72%  '
$ws.Range("D45").Value = '''0.0577'
$ws.Range("E45").Value = '  +8.78%  '
$ws.Range("D46").Value = '''0.101'
$ws.Range("E46").Value = '  +0.34%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").Value = '''1.00'
$ws.Range("E47").Value = '  +0.49%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''20.04'
$ws.Range("E48").Value = '  +15.25%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''5.03'
$ws.Range("E49").Value = '  +7.65%  '
$ws.Range("E50").Value = '  +4.24%  '
$ws.Range("D51").Value = '''2.039.06'
$ws.Range("E51").Value = '  +9.77%  '
